$d = $word.ActiveDocument

# Shared pkg:package wrapper used by every InsertXML call below. InsertXML
# replaces the *entire* contents of the target Range, so each payload is a
# complete replacement <w:p> (including its original w14:paraId / rsidR
# bookkeeping attributes and pPr) built from the paragraph's current runs,
# minus the <w:proofErr/> markers we want gone.
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Locate the "...2022.06.09 周四" paragraph that still carries the
# proofErr wrapper (the first / earlier of the two identical-looking
# paragraphs in the document) and the "...2022.06.17 周五" paragraph.
$paraZhouSi = $null
$paraZhouWu = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($paraZhouSi -eq $null -and $t -like "*2022.06.09*" -and $t -like "*周四*") {
        $paraZhouSi = $i
    }
    if ($paraZhouWu -eq $null -and $t -like "*2022.06.1*" -and $t -like "*周五*") {
        $paraZhouWu = $i
    }
}

# --- Edit 1: drop the spellStart/spellEnd proofErr wrapper around "周四";
#     everything else in the paragraph (including the date) is unchanged.
$p1 = $d.Paragraphs($paraZhouSi)
$r1 = $p1.Range
$xml1 = $pkgOpen + '<w:p w14:paraId="40B1E3BF" w14:textId="77777777" w:rsidR="00EE5398" w:rsidRDefault="00395A0A"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>预计完成时间：</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">2022.06.09 </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>周四</w:t></w:r></w:p>' + $pkgClose
$r1.InsertXML($xml1)

# --- Edit 2: date "2022.06.17" -> "2022.06.24" (keeping the existing
#     "2022.06." / day-digits run split, and the w:rsidR="00BF5DEE" stamp
#     on the run carrying the day digits), and drop the proofErr wrapper
#     around "周五".
$p2 = $d.Paragraphs($paraZhouWu)
$r2 = $p2.Range
$xml2 = $pkgOpen + '<w:p w14:paraId="49BA7EE9" w14:textId="71DD6C78" w:rsidR="00EE5398" w:rsidRDefault="00395A0A"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>预计完成时间：</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2022.06.</w:t></w:r><w:r w:rsidR="00BF5DEE"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>24</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>周五</w:t></w:r></w:p>' + $pkgClose
$r2.InsertXML($xml2)

Write-Output ("edited paragraphs " + $paraZhouSi + " and " + $paraZhouWu)
